# Auto-applies the crypto price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.239.16'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.425.76'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.57%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.40'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.79'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.37%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("E9").Value = '  +5.64%  '
$ws.Range("E10").Value = '  +0.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '48.24'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.14%  '
$ws.Range("E12").Value = '  +2.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '678.34'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.62%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.978.44'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.60%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.66'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.369.48'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.91%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.422.01'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.09%  '
$ws.Range("E18").Value = '  +0.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.79'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.33'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.912'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.38'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.01'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '100.75'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.25%  '
$ws.Range("E25").Value = '  -0.56%  '
$ws.Range("E26").Value = '  -0.82%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.68'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.83%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '33.62'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.77'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.77%  '
$ws.Range("E30").Value = '  -1.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.74'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +9.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '560.03'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.19%  '
$ws.Range("E33").Value = '  -1.10%  '
$ws.Range("E34").Value = '  -0.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.03'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.67%  '
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.619.10'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.82%  '
$ws.Range("E38").Value = '  -2.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.03'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0733'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.06%  '
$ws.Range("E41").Value = '  +2.13%  '
$ws.Range("E42").Value = '  +1.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.40'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.26%  '
$ws.Range("E44").Value = '  +1.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.335'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.76%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.67'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.37%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.41'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.90%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.129'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.08%  '
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '130.97'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.68'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.83%  '
